$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update BOM row for the combined connector designators (J16, J15 added) ---
# Designator list (B2): J10,J12,J11  ->  J16,J15,J10,J12,J11
$ws.Range("B2").Value = "J16,J15,J10,J12,J11"

# Quantity (A2) now reflects the 5 designators instead of 3
$ws.Range("A2").Value = 5

# The designator cell picked up an explicit "Noto Sans" font while being edited
$ws.Range("B2").Font.Name = "Noto Sans"

# --- Preserve page setup details that were explicitly set on the sheet ---
$ps = $ws.PageSetup
$ps.PaperSize = 1
$ps.Zoom = 100
$ps.FitToPagesWide = 1
$ps.FitToPagesTall = 1
$ps.Orientation = 1

# --- Selection moved to C18 when the author saved the file ---
$ws.Range("C18").Select()
